$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 39.8
$ws.Range("I5").Value = 39.8
$ws.Range("K5").Value = 39.8
$ws.Range("M5").Value = 75.2
$ws.Range("H43").Value = 2203
$ws.Range("I43").Value = 1730
$ws.Range("K43").Value = 1730
$ws.Range("M43").Value = -1661
$ws.Range("H51").Value = 2000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -2968
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H74").Value = 6938.6333
$ws.Range("I74").Value = 4568.9287
$ws.Range("K74").Value = 4568.9287
$ws.Range("M74").Value = -3632.9287
$ws.Range("H77").Value = 6938.6333
$ws.Range("I77").Value = 4568.9287
$ws.Range("K77").Value = 22844.6435
$ws.Range("M77").Value = -18164.6435
$ws.Range("H100").Value = 5969.0557
$ws.Range("J100").Value = 9577.875
$ws.Range("L100").Value = 9577.875
$ws.Range("N100").Value = -10659.875
$ws.Range("H112").Value = 1182.591
$ws.Range("J112").Value = 1200.8536
$ws.Range("L112").Value = 3602.5608
$ws.Range("N112").Value = -5818.560799999999
$ws.Range("H132").Value = 4165.7607
$ws.Range("I132").Value = 3060.4524
$ws.Range("K132").Value = 9181.3572
$ws.Range("M132").Value = -6651.3572
$ws.Range("H135").Value = 1087.5625
$ws.Range("I135").Value = 1038.4166
$ws.Range("J135").Value = 1235
$ws.Range("K135").Value = 9345.749400000001
$ws.Range("L135").Value = 11115
$ws.Range("M135").Value = -6810.749400000001
$ws.Range("N135").Value = -16185
$ws.Range("H138").Value = 2716.5476
$ws.Range("I138").Value = 1467
$ws.Range("J138").Value = 3410.7407
$ws.Range("K138").Value = 4401
$ws.Range("L138").Value = 10232.2221
$ws.Range("M138").Value = 739
$ws.Range("N138").Value = -20512.2221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1066.3334
$ws.Range("I12").Value = 1099.5
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 1099.5
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -926.5
$ws.Range("N12").Value = -1346
$ws.Range("H74").Value = 2472.2964
$ws.Range("I74").Value = 1341.55
$ws.Range("K74").Value = 1341.55
$ws.Range("M74").Value = -467.55
$ws.Range("H77").Value = 2472.2964
$ws.Range("I77").Value = 1341.55
$ws.Range("K77").Value = 6707.75
$ws.Range("M77").Value = -2339.75
$ws.Range("H102").Value = 2756.8572
$ws.Range("I102").Value = 2618
$ws.Range("J102").Value = 3266
$ws.Range("K102").Value = 2618
$ws.Range("L102").Value = 3266
$ws.Range("M102").Value = -996
$ws.Range("N102").Value = -6510
$ws.Range("H132").Value = 2390.5454
$ws.Range("I132").Value = 1746.4706
$ws.Range("J132").Value = 4580.4
$ws.Range("K132").Value = 5239.4118
$ws.Range("L132").Value = 13741.2
$ws.Range("M132").Value = -2709.4118
$ws.Range("N132").Value = -18801.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 700
$ws.Range("J7").Value = 700
$ws.Range("L7").Value = 700
$ws.Range("N7").Value = -926
$ws.Range("H99").Value = 9246.579
$ws.Range("I99").Value = 4621.231
$ws.Range("J99").Value = 19268.166
$ws.Range("K99").Value = 4621.231
$ws.Range("L99").Value = 19268.166
$ws.Range("M99").Value = -3123.231
$ws.Range("N99").Value = -22264.166
$ws.Range("H107").Value = 3655.95
$ws.Range("I107").Value = 3595.2354
$ws.Range("K107").Value = 3595.2354
$ws.Range("M107").Value = -1675.2354
$ws.Range("H134").Value = 2332.55
$ws.Range("I134").Value = 1447.0303
$ws.Range("K134").Value = 4341.090899999999
$ws.Range("M134").Value = -1806.090899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 633.3333
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H31").Value = 4739.84
$ws.Range("I31").Value = 2077.8333
$ws.Range("J31").Value = 7197.077
$ws.Range("K31").Value = 2077.8333
$ws.Range("L31").Value = 7197.077
$ws.Range("M31").Value = -1782.8333
$ws.Range("N31").Value = -7787.077
$ws.Range("H34").Value = 4739.84
$ws.Range("I34").Value = 2077.8333
$ws.Range("J34").Value = 7197.077
$ws.Range("K34").Value = 2077.8333
$ws.Range("L34").Value = 7197.077
$ws.Range("M34").Value = -1875.8333
$ws.Range("N34").Value = -7601.077
$ws.Range("H58").Value = 3874.95
$ws.Range("I58").Value = 1313.8334
$ws.Range("J58").Value = 7716.625
$ws.Range("K58").Value = 1313.8334
$ws.Range("L58").Value = 7716.625
$ws.Range("M58").Value = -1110.8334
$ws.Range("N58").Value = -8122.625
$ws.Range("H134").Value = 2268
$ws.Range("I134").Value = 2086.1304
$ws.Range("K134").Value = 6258.3912
$ws.Range("M134").Value = -3723.3912
$ws.Range("H136").Value = 3874.95
$ws.Range("I136").Value = 1313.8334
$ws.Range("J136").Value = 7716.625
$ws.Range("K136").Value = 3941.5002
$ws.Range("L136").Value = 23149.875
$ws.Range("M136").Value = -1391.5002
$ws.Range("N136").Value = -28249.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10653.857
$ws.Range("I3").Value = 1315.4
$ws.Range("K3").Value = 3946.2
$ws.Range("M3").Value = -3834.2
$ws.Range("H107").Value = 230.25
$ws.Range("I107").Value = 220.28572
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 660.85716
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1259.14284
$ws.Range("N107").Value = -4740
$ws.Range("H121").Value = 22223154
$ws.Range("I121").Value = 50000164
$ws.Range("J121").Value = 1543.8
$ws.Range("K121").Value = 150000492
$ws.Range("L121").Value = 4631.4
$ws.Range("M121").Value = -149999182
$ws.Range("N121").Value = -7251.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 202.61111
$ws.Range("I2").Value = 161.76923
$ws.Range("J2").Value = 308.8
$ws.Range("K2").Value = 161.76923
$ws.Range("L2").Value = 308.8
$ws.Range("M2").Value = -48.76922999999999
$ws.Range("N2").Value = -534.8
$ws.Range("H9").Value = 6906.1
$ws.Range("J9").Value = 20966.334
$ws.Range("L9").Value = 20966.334
$ws.Range("N9").Value = -21306.334
$ws.Range("H70").Value = 7128.4287
$ws.Range("I70").Value = 6974.75
$ws.Range("K70").Value = 6974.75
$ws.Range("M70").Value = -6704.75
$ws.Range("H73").Value = 7128.4287
$ws.Range("I73").Value = 6974.75
$ws.Range("K73").Value = 6974.75
$ws.Range("M73").Value = -6038.75
$ws.Range("H80").Value = 377534.3
$ws.Range("I80").Value = 627256.4
$ws.Range("J80").Value = 127812.25
$ws.Range("K80").Value = 627256.4
$ws.Range("L80").Value = 127812.25
$ws.Range("M80").Value = -626258.4
$ws.Range("N80").Value = -129808.25
$ws.Range("H83").Value = 377534.3
$ws.Range("I83").Value = 627256.4
$ws.Range("J83").Value = 127812.25
$ws.Range("K83").Value = 3136282
$ws.Range("L83").Value = 639061.25
$ws.Range("M83").Value = -3131290
$ws.Range("N83").Value = -649045.25
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H131").Value = 60308.668
$ws.Range("J131").Value = 60308.668
$ws.Range("L131").Value = 60308.668
$ws.Range("N131").Value = -70388.66800000001
$ws.Range("H132").Value = 13168923
$ws.Range("I132").Value = 17252270
$ws.Range("J132").Value = 11475.389
$ws.Range("K132").Value = 51756810
$ws.Range("L132").Value = 34426.167
$ws.Range("M132").Value = -51754280
$ws.Range("N132").Value = -39486.167
$ws.Range("H136").Value = 49326
$ws.Range("J136").Value = 49326
$ws.Range("L136").Value = 147978
$ws.Range("N136").Value = -153078

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5836.579
$ws.Range("I16").Value = 9753.091
$ws.Range("K16").Value = 9753.091
$ws.Range("M16").Value = -9583.091
$ws.Range("H68").Value = 3085.7144
$ws.Range("J68").Value = 3400
$ws.Range("L68").Value = 3400
$ws.Range("N68").Value = -4898
$ws.Range("H71").Value = 3085.7144
$ws.Range("J71").Value = 3400
$ws.Range("L71").Value = 17000
$ws.Range("N71").Value = -24488
$ws.Range("H82").Value = 2550.3076
$ws.Range("I82").Value = 1877.4
$ws.Range("K82").Value = 1877.4
$ws.Range("M82").Value = -1516.4
$ws.Range("H85").Value = 2550.3076
$ws.Range("I85").Value = 1877.4
$ws.Range("K85").Value = 1877.4
$ws.Range("M85").Value = -629.4000000000001
$ws.Range("H93").Value = 2486.5
$ws.Range("I93").Value = 2583.8
$ws.Range("K93").Value = 2583.8
$ws.Range("M93").Value = -1335.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3262.7368
$ws.Range("I107").Value = 1606.8462
$ws.Range("K107").Value = 4820.5386
$ws.Range("M107").Value = -2900.5386
$ws.Range("H123").Value = 77000
$ws.Range("J123").Value = 77000
$ws.Range("L123").Value = 77000
$ws.Range("N123").Value = -86800
$ws.Range("H132").Value = 3451244.5
$ws.Range("I132").Value = 3923728.8
$ws.Range("K132").Value = 11771186.4
$ws.Range("M132").Value = -11768656.4
